$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G6").Value = 3.1
$ws.Range("H6").Value = 2.75
$ws.Range("I6").Value = 2.63
$ws.Range("U6").Value = 2.38
$ws.Range("V6").Value = 1.53
$ws.Range("AC6").Value = 5
$ws.Range("AJ6").Value = 26
$ws.Range("AP6").Value = 41
$ws.Range("AU6").Value = 10
$ws.Range("AW6").Value = 4.33
